$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.158.28"
$ws.Range("E2").Value = "  +2.92%  "
$ws.Range("D3").Value = "2.277.39"
$ws.Range("E3").Value = "  +3.01%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.32%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "319.63"
$ws.Range("E5").Value = "  +1.47%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "103.88"
$ws.Range("E6").Value = "  +6.46%  "
$ws.Range("E7").Value = "  +1.48%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  -0.21%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.573"
$ws.Range("E9").Value = "  +3.04%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.88"
$ws.Range("E10").Value = "  +6.71%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0840"
$ws.Range("E11").Value = "  +1.66%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.91"
$ws.Range("E12").Value = "  +2.35%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.107"
$ws.Range("E13").Value = "  +2.11%  "
$ws.Range("D14").Value = "2.624.00"
$ws.Range("E14").Value = "  +2.98%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.879"
$ws.Range("E15").Value = "  +2.03%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.56"
$ws.Range("E16").Value = "  +3.56%  "
$ws.Range("D17").Value = "2.277.83"
$ws.Range("E17").Value = "  +3.37%  "
$ws.Range("D18").Value = "44.196.44"
$ws.Range("E18").Value = "  +3.07%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.18"
$ws.Range("E19").Value = "  -2.47%  "
$ws.Range("D20").Value = "0.0₂01000"
$ws.Range("E20").Value = "  +4.55%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.62"
$ws.Range("E21").Value = "  +3.90%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "66.29"
$ws.Range("E22").Value = "  +1.64%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.22"
$ws.Range("E23").Value = "  +1.55%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "237.83"
$ws.Range("E24").Value = "  +0.73%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.20"
$ws.Range("E25").Value = "  +4.10%  "
$ws.Range("E26").Value = "  +0.16%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.28"
$ws.Range("E27").Value = "  +2.72%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "39.09"
$ws.Range("E28").Value = "  +16.17%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.21"
$ws.Range("E29").Value = "  -0.32%  "
$ws.Range("B30").Value = "Filecoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.56"
$ws.Range("E30").Value = "  +5.81%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "162.46"
$ws.Range("E31").Value = "  +5.16%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.52"
$ws.Range("E32").Value = "  +0.58%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0883"
$ws.Range("E33").Value = "  +0.91%  "
$ws.Range("E34").Value = "  -2.26%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.08"
$ws.Range("E35").Value = "  +5.04%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.26"
$ws.Range("E36").Value = "  +2.17%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.111"
$ws.Range("E37").Value = "  +8.90%  "
$ws.Range("E38").Value = "  -0.22%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.53"
$ws.Range("E39").Value = "  +2.84%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.91"
$ws.Range("E40").Value = "  +6.14%  "
$ws.Range("E41").Value = "  +1.18%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "15.44"
$ws.Range("E42").Value = "  +26.24%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.00"
$ws.Range("E43").Value = "  -0.13%  "
$ws.Range("D44").Value = "1.778.54"
$ws.Range("E44").Value = "  -4.46%  "
$ws.Range("E45").Value = "  +0.88%  "
$ws.Range("B46").Value = "BitcoinSV"
$ws.Range("C46").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "85.28"
$ws.Range("E46").Value = "  -3.77%  "
$ws.Range("B47").Value = "THORChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "5.37"
$ws.Range("E47").Value = "  -0.64%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.89"
$ws.Range("E48").Value = "  +3.20%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "59.97"
$ws.Range("E49").Value = "  +0.21%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "75.10"
$ws.Range("E50").Value = "  -1.06%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "104.65"
$ws.Range("E51").Value = "  +4.17%  "
